# 9th Stab - Cosmetic Changes
#
# The "2018" sheet currently has three used columns:
#   A = analyst name, B = rating ("UN"), C = date/action details
# with a one-row header in B1/C1 holding date labels ("Jun_13", "Jun_10").
#
# This change inserts two more "UN" rating columns (one per extra
# observation date) immediately before the existing rating/detail column,
# shifting the old column C to column E and adding fresh header labels
# "Jun_15" / "Jun_17" plus "UN" placeholder values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C; this pushes the existing column C
# (dates/details, including its per-row fill styling) out to column E.
$ws.Columns.Item(3).EntireColumn.Insert()
$ws.Columns.Item(3).EntireColumn.Insert()

# Give the two new columns the same raw column width as the rest of the
# table (stored width 8.0 "characters") so they render consistently.
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667

# Header row: new "Jun_15"/"Jun_17" labels slot in ahead of the existing
# "Jun_13"/"Jun_10" headers (which simply shift right with the columns).
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Data rows: fill the two newly inserted columns with the same default
# "UN" placeholder used throughout column B.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
